$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: Create the new "2022-Q3" sheet by duplicating the existing
# "2022-Q2" sheet (same column layout/styling) and inserting the copy
# directly before it, then trim it down and overwrite its values.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)           # "总计"
$q2Sheet = $wb.Worksheets.Item(2)           # "2022-Q2"

$q2Sheet.Copy($q2Sheet)                     # duplicate inserted before "2022-Q2"
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# "2022-Q2" has 9 data rows (rows 2-10); "2022-Q3" only needs 4 (rows 2-5).
$q3Sheet.Range("A6:H10").EntireRow.Delete()

# Row 2: 159792 富国中证港股通互联网ETF
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "'159792"
$q3Sheet.Range("B2").ClearFormats()
$q3Sheet.Range("C2").Value = "富国中证港股通互联网ETF"
$q3Sheet.Range("D2").Value = "'18.08"
$q3Sheet.Range("D2").ClearFormats()
$q3Sheet.Range("E2").Value = "'99.26"
$q3Sheet.Range("E2").ClearFormats()
$q3Sheet.Range("F2").Value = "'3.95"
$q3Sheet.Range("F2").ClearFormats()
$q3Sheet.Range("G2").Value = "'0.7142"
$q3Sheet.Range("G2").ClearFormats()
$q3Sheet.Range("H2").Value = 8

# Row 3: 513770 华宝中证港股通互联网ETF
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "'513770"
$q3Sheet.Range("B3").ClearFormats()
$q3Sheet.Range("C3").Value = "华宝中证港股通互联网ETF"
$q3Sheet.Range("D3").Value = "'3.80"
$q3Sheet.Range("D3").ClearFormats()
$q3Sheet.Range("E3").Value = "'98.21"
$q3Sheet.Range("E3").ClearFormats()
$q3Sheet.Range("F3").Value = "'3.86"
$q3Sheet.Range("F3").ClearFormats()
$q3Sheet.Range("G3").Value = "'0.1467"
$q3Sheet.Range("G3").ClearFormats()
$q3Sheet.Range("H3").Value = 8

# Row 4: 003993 前海开源沪港深核心驱动灵活配置混合
$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = "'003993"
$q3Sheet.Range("B4").ClearFormats()
$q3Sheet.Range("C4").Value = "前海开源沪港深核心驱动灵活配置混合"
$q3Sheet.Range("D4").Value = "'0.53"
$q3Sheet.Range("D4").ClearFormats()
$q3Sheet.Range("E4").Value = "'82.41"
$q3Sheet.Range("E4").ClearFormats()
$q3Sheet.Range("F4").Value = "'7.14"
$q3Sheet.Range("F4").ClearFormats()
$q3Sheet.Range("G4").Value = "'0.0378"
$q3Sheet.Range("G4").ClearFormats()
$q3Sheet.Range("H4").Value = 4

# Row 5: 003413 华泰柏瑞新经济沪港深混合
$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = "'003413"
$q3Sheet.Range("B5").ClearFormats()
$q3Sheet.Range("C5").Value = "华泰柏瑞新经济沪港深混合"
$q3Sheet.Range("D5").Value = "'0.42"
$q3Sheet.Range("D5").ClearFormats()
$q3Sheet.Range("E5").Value = "'86.45"
$q3Sheet.Range("E5").ClearFormats()
$q3Sheet.Range("F5").Value = "'5.48"
$q3Sheet.Range("F5").ClearFormats()
$q3Sheet.Range("G5").Value = "'0.0230"
$q3Sheet.Range("G5").ClearFormats()
$q3Sheet.Range("H5").Value = 6

# ------------------------------------------------------------------
# Step 2: Update the "总计" (summary) sheet with the new quarter's
# row and shift the existing history down by one row.
# ------------------------------------------------------------------
# Prepare the brand-new row 9 (previously nonexistent) so that column
# A keeps the same bold/bordered style used by the other index cells.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.92

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 9
$summary.Range("D3").Value = 4.01

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 6
$summary.Range("D4").Value = 0.92

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 0.36

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.21

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 1
$summary.Range("D7").Value = 0.05

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 6
$summary.Range("D8").Value = 1.95

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 3
$summary.Range("D9").Value = 1.4
